# Applies the LOM3248 content update:
#  - "Ativação:" date changes from 01/01/2012 to 01/01/2023 (rows 8 and 13,
#    which share the same underlying text).
#  - Three new English paragraphs are added alongside their section headers
#    in column A (rows 11, 14 and 16), mirrored into both column B (the
#    "modified" value) and column C (the same value shown in red).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (Objectives / Objectives:) : add English objectives text -------
# Copy formatting from an existing populated row in the same columns first,
# so the new cells pick up the normal (B) / red (C) wrapped-text styles.
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."
$ws.Range("C11").Value = "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."

# --- Row 14 (Short syllabus:) : add short syllabus text ---------------------
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B14").Value = "To be defined, according to the programmed topic."
$ws.Range("C14").Value = "To be defined, according to the programmed topic."

# --- Row 16 (Syllabus:) : add full syllabus text -----------------------------
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("B16").Value = "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Range("C16").Value = "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."

# --- Update the activation date (shared by rows 8 and 13) ------------------
# Force the cells to remain plain text (not auto-converted to a date
# serial) by formatting them as Text before assigning the new value. Done
# last so it doesn't leak its "Text" number format into the copies above.
$dateRange = $ws.Range("B8:C8")
$dateRange.NumberFormat = "@"
$dateRange.Value = "01/01/2023"

$dateRange13 = $ws.Range("B13:C13")
$dateRange13.NumberFormat = "@"
$dateRange13.Value = "01/01/2023"

$excel.CutCopyMode = 0
